$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells whose new numeric-looking values
# must remain text strings (matches original inlineStr text cells)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "34.179.83"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "1.781.59"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "226.27"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "31.82"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  +1.43%  "
$ws.Range("D12").Value = "2.038.16"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "11.02"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.788.72"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "34.132.52"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").Value = "247.11"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "11.00"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "162.71"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "16.32"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "0.0520"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("E33").Value = "  +5.03%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "1.444.01"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("D36").Value = "0.654"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  +7.77%  "
$ws.Range("E38").Value = "  +3.46%  "
$ws.Range("D39").Value = "1.04"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "80.28"
$ws.Range("E41").Value = "  +2.44%  "
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "13.67"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").Value = "6.07"
$ws.Range("E46").Value = "  +3.94%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "0.0₆0137"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").Value = "1.940.96"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  +0.22%  "
